$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (rows 2-51) to remain plain text so numeric-looking
# strings (e.g. "1.00", "10.80") are not coerced into numbers, matching
# the original inline-string cell typing used throughout this sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "69.747.06"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "2.444.84"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "565.98"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "166.58"
$ws.Range("E6").Value = "  +2.20%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.175"
$ws.Range("E9").Value = "  +10.22%  "

$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "2.445.04"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.162"
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").Value = "4.68"
$ws.Range("E13").Value = "  -2.28%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000182"
$ws.Range("E14").Value = "  +6.93%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "69.639.54"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.890.92"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "24.04"
$ws.Range("E17").Value = "  +3.77%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.442.01"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "10.80"
$ws.Range("E19").Value = "  +4.11%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "342.16"
$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "7.12"
$ws.Range("E21").Value = "  +3.65%  "

$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").Value = "3.88"
$ws.Range("E22").Value = "  +2.30%  "

$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  +5.86%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "66.26"
$ws.Range("E25").Value = "  -0.84%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "3.88"
$ws.Range("E26").Value = "  +5.43%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "8.53"
$ws.Range("E27").Value = "  +5.15%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.568.99"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.02"
$ws.Range("E29").Value = "  +2.41%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0858"
$ws.Range("E30").Value = "  +5.51%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "7.33"
$ws.Range("E31").Value = "  +2.07%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  +9.60%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "449.91"
$ws.Range("E33").Value = "  +5.80%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "162.08"
$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D37").Value = "19.06"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.109"
$ws.Range("E39").Value = "  +5.42%  "

$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.14"
$ws.Range("E40").Value = "  +2.16%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "0.304"
$ws.Range("E41").Value = "  +3.06%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.54"
$ws.Range("E42").Value = "  +4.97%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "4.44"
$ws.Range("E43").Value = "  +2.47%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").Value = "1.09"
$ws.Range("E44").Value = "  +3.50%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.15"
$ws.Range("E45").Value = "  +7.15%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "132.30"
$ws.Range("E47").Value = "  +2.20%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0725"
$ws.Range("E48").Value = "  +1.30%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.489"
$ws.Range("E49").Value = "  +2.37%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.561"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.0932"
$ws.Range("E51").Value = "  +1.65%  "

